$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 15.62169031819198
$ws.Cells.Item(2, 3).Value = 12.59394899082903
$ws.Cells.Item(2, 5).Value = 9.62753802310408
$ws.Cells.Item(2, 6).Value = 16.86991607391245
$ws.Cells.Item(2, 7).Value = 3.639338297477192
$ws.Cells.Item(2, 9).Value = 23.69129317225615
$ws.Cells.Item(2, 13).Value = 16.48787687408263
$ws.Cells.Item(2, 14).Value = 17.96594314222329
$ws.Cells.Item(3, 2).Value = 15.02176847277573
$ws.Cells.Item(3, 3).Value = 11.93515176590816
$ws.Cells.Item(3, 5).Value = 9.540505812403776
$ws.Cells.Item(3, 6).Value = 15.89584955866815
$ws.Cells.Item(3, 7).Value = 3.643554000904162
$ws.Cells.Item(3, 9).Value = 23.63126395042023
$ws.Cells.Item(3, 13).Value = 16.1988686286642
$ws.Cells.Item(3, 14).Value = 18.03568555176137
$ws.Cells.Item(4, 2).Value = 14.64605423773752
$ws.Cells.Item(4, 3).Value = 11.51559996041644
$ws.Cells.Item(4, 5).Value = 9.489677000107603
$ws.Cells.Item(4, 6).Value = 15.26997757108489
$ws.Cells.Item(4, 7).Value = 3.646272119687463
$ws.Cells.Item(4, 9).Value = 23.60164956058836
$ws.Cells.Item(4, 13).Value = 16.02424758217619
$ws.Cells.Item(4, 14).Value = 18.08051894385994
$ws.Cells.Item(5, 2).Value = 14.49136409714109
$ws.Cells.Item(5, 3).Value = 11.34107791868113
$ws.Cells.Item(5, 5).Value = 9.469637546380721
$ws.Cells.Item(5, 6).Value = 15.008197319934
$ws.Cells.Item(5, 7).Value = 3.647412518505437
$ws.Cells.Item(5, 9).Value = 23.59140080532088
$ws.Cells.Item(5, 13).Value = 15.95389819029981
$ws.Cells.Item(5, 14).Value = 18.09929518207692
$ws.Cells.Item(6, 2).Value = 14.46559058842338
$ws.Cells.Item(6, 3).Value = 11.3118917508655
$ws.Cells.Item(6, 5).Value = 9.466351183665045
$ws.Cells.Item(6, 6).Value = 14.96433081551589
$ws.Cells.Item(6, 7).Value = 3.647603862561392
$ws.Cells.Item(6, 9).Value = 23.58980878528501
$ws.Cells.Item(6, 13).Value = 15.94226850364396
$ws.Cells.Item(6, 14).Value = 18.10244355506729
$ws.Cells.Item(7, 2).Value = 14.6439740720692
$ws.Cells.Item(7, 3).Value = 11.51326033793285
$ws.Cells.Item(7, 5).Value = 9.489403991138836
$ws.Cells.Item(7, 6).Value = 15.26647399323726
$ws.Cells.Item(7, 7).Value = 3.646287366742223
$ws.Cells.Item(7, 9).Value = 23.60150398024161
$ws.Cells.Item(7, 13).Value = 16.02329541661462
$ws.Cells.Item(7, 14).Value = 18.08077011619107
$ws.Cells.Item(8, 2).Value = 15.41651706410171
$ws.Cells.Item(8, 3).Value = 12.37004662727982
$ws.Cells.Item(8, 5).Value = 9.596998744276943
$ws.Cells.Item(8, 6).Value = 16.53996406344765
$ws.Cells.Item(8, 7).Value = 3.640765045233453
$ws.Cells.Item(8, 9).Value = 23.66909099573003
$ws.Cells.Item(8, 13).Value = 16.38769836366198
$ws.Cells.Item(8, 14).Value = 17.98957319549417
$ws.Cells.Item(9, 2).Value = 16.8626490080558
$ws.Cells.Item(9, 3).Value = 13.92242691862002
$ws.Cells.Item(9, 5).Value = 9.82785477807516
$ws.Cells.Item(9, 6).Value = 19.0027458068253
$ws.Cells.Item(9, 7).Value = 3.630958212891828
$ws.Cells.Item(9, 9).Value = 23.85910459670133
$ws.Cells.Item(9, 13).Value = 17.12042194314
$ws.Cells.Item(9, 14).Value = 17.82667013403235
$ws.Cells.Item(10, 2).Value = 17.87098504893731
$ws.Cells.Item(10, 3).Value = 14.9760436366061
$ws.Cells.Item(10, 5).Value = 10.00836224012229
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 3.624367494676913
$ws.Cells.Item(10, 9).Value = 24.03360754789213
$ws.Cells.Item(10, 13).Value = 17.66405249523809
$ws.Cells.Item(10, 14).Value = 17.7166629324714
$ws.Cells.Item(11, 2).Value = 18.31585227159218
$ws.Cells.Item(11, 3).Value = 15.43510796548865
$ws.Cells.Item(11, 5).Value = 10.09256754473292
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 3.621500690965345
$ws.Cells.Item(11, 9).Value = 24.12049616638838
$ws.Cells.Item(11, 13).Value = 17.91133220703642
$ws.Cells.Item(11, 14).Value = 17.66871218994819
$ws.Cells.Item(12, 2).Value = 18.48217773188955
$ws.Cells.Item(12, 3).Value = 15.60594697723754
$ws.Cells.Item(12, 5).Value = 10.12473021452458
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 3.620433845925582
$ws.Cells.Item(12, 9).Value = 24.15446805531381
$ws.Cells.Item(12, 13).Value = 18.00487311936127
$ws.Cells.Item(12, 14).Value = 17.65085497878833
$ws.Cells.Item(13, 2).Value = 18.4464537108219
$ws.Cells.Item(13, 3).Value = 15.56928834910186
$ws.Cells.Item(13, 5).Value = 10.11779153500606
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 3.62066277829235
$ws.Cells.Item(13, 9).Value = 24.1471042417087
$ws.Cells.Item(13, 13).Value = 17.98473327330362
$ws.Cells.Item(13, 14).Value = 17.65468747700918
$ws.Cells.Item(14, 2).Value = 18.32957955252415
$ws.Cells.Item(14, 3).Value = 15.44922354313729
$ws.Cells.Item(14, 5).Value = 10.09520820022357
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 3.621412545957985
$ws.Cells.Item(14, 9).Value = 24.12326969543458
$ws.Cells.Item(14, 13).Value = 17.91903038282107
$ws.Cells.Item(14, 14).Value = 17.66723704074014
$ws.Cells.Item(15, 2).Value = 18.25770853388747
$ws.Cells.Item(15, 3).Value = 15.37528750592183
$ws.Cells.Item(15, 5).Value = 10.08141044805311
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 3.621874238233836
$ws.Cells.Item(15, 9).Value = 24.10880922707011
$ws.Cells.Item(15, 13).Value = 17.87876986102504
$ws.Cells.Item(15, 14).Value = 17.67496316667962
$ws.Cells.Item(16, 2).Value = 17.84162058445389
$ws.Cells.Item(16, 3).Value = 14.94562754827087
$ws.Cells.Item(16, 5).Value = 10.00289920259239
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 3.624557478432792
$ws.Cells.Item(16, 9).Value = 24.0280793828675
$ws.Cells.Item(16, 13).Value = 17.64788351194566
$ws.Cells.Item(16, 14).Value = 17.7198387108644
$ws.Cells.Item(17, 2).Value = 17.58271070347613
$ws.Cells.Item(17, 3).Value = 14.67679402321719
$ws.Cells.Item(17, 5).Value = 9.955253297751357
$ws.Cells.Item(17, 6).Value = 20.20408069617459
$ws.Cells.Item(17, 7).Value = 3.626237103555924
$ws.Cells.Item(17, 9).Value = 23.98047010516683
$ws.Cells.Item(17, 13).Value = 17.50616774024779
$ws.Cells.Item(17, 14).Value = 17.74790411837697
$ws.Cells.Item(18, 2).Value = 17.43249833450491
$ws.Cells.Item(18, 3).Value = 14.52026795214297
$ws.Cells.Item(18, 5).Value = 9.928047022972002
$ws.Cells.Item(18, 6).Value = 19.95656407809808
$ws.Cells.Item(18, 7).Value = 3.627215550652652
$ws.Cells.Item(18, 9).Value = 23.9537937979677
$ws.Cells.Item(18, 13).Value = 17.42466013684267
$ws.Cells.Item(18, 14).Value = 17.76424352121046
$ws.Cells.Item(19, 2).Value = 17.3814215360375
$ws.Cells.Item(19, 3).Value = 14.46694758147562
$ws.Cells.Item(19, 5).Value = 9.918870267420141
$ws.Cells.Item(19, 6).Value = 19.87204792380562
$ws.Cells.Item(19, 7).Value = 3.627548964935105
$ws.Cells.Item(19, 9).Value = 23.94488341695276
$ws.Cells.Item(19, 13).Value = 17.39706664923523
$ws.Cells.Item(19, 14).Value = 17.7698095890717
$ws.Cells.Item(20, 2).Value = 17.61040718304453
$ws.Cells.Item(20, 3).Value = 14.70560917162223
$ws.Cells.Item(20, 5).Value = 9.960304931359424
$ws.Cells.Item(20, 6).Value = 20.2495528364879
$ws.Cells.Item(20, 7).Value = 3.626057025214379
$ws.Cells.Item(20, 9).Value = 23.98546506049783
$ws.Cells.Item(20, 13).Value = 17.52125396231816
$ws.Cells.Item(20, 14).Value = 17.74489612829501
$ws.Cells.Item(21, 2).Value = 18.36396737650603
$ws.Cells.Item(21, 3).Value = 15.48457149812819
$ws.Cells.Item(21, 5).Value = 10.10183418387652
$ws.Cells.Item(21, 6).Value = 21.46857628470567
$ws.Cells.Item(21, 7).Value = 3.621191813200313
$ws.Cells.Item(21, 9).Value = 24.13024155845835
$ws.Cells.Item(21, 13).Value = 17.93833233133569
$ws.Cells.Item(21, 14).Value = 17.66354276981352
$ws.Cells.Item(22, 2).Value = 18.84395677884638
$ws.Cells.Item(22, 3).Value = 15.97615988710617
$ws.Cells.Item(22, 5).Value = 10.19592838863483
$ws.Cells.Item(22, 6).Value = 22.22866616901555
$ws.Cells.Item(22, 7).Value = 3.618121349676367
$ws.Cells.Item(22, 9).Value = 24.23108637616426
$ws.Cells.Item(22, 13).Value = 18.21030127448633
$ws.Cells.Item(22, 14).Value = 17.61212629172613
$ws.Cells.Item(23, 2).Value = 18.58896538599786
$ws.Cells.Item(23, 3).Value = 15.71541716344802
$ws.Cells.Item(23, 5).Value = 10.14557074975465
$ws.Cells.Item(23, 6).Value = 21.82633154475864
$ws.Cells.Item(23, 7).Value = 3.619750163947109
$ws.Cells.Item(23, 9).Value = 24.17669796714421
$ws.Cells.Item(23, 13).Value = 18.06523322821621
$ws.Cells.Item(23, 14).Value = 17.63940790170263
$ws.Cells.Item(24, 2).Value = 17.59788983441378
$ws.Cells.Item(24, 3).Value = 14.69258797247464
$ws.Cells.Item(24, 5).Value = 9.95802050660906
$ws.Cells.Item(24, 6).Value = 20.22900810905294
$ws.Cells.Item(24, 7).Value = 3.626138398727343
$ws.Cells.Item(24, 9).Value = 23.98320467597248
$ws.Cells.Item(24, 13).Value = 17.51443357780439
$ws.Cells.Item(24, 14).Value = 17.74625540448077
$ws.Cells.Item(25, 2).Value = 16.48015252221285
$ws.Cells.Item(25, 3).Value = 13.51718590269832
$ws.Cells.Item(25, 5).Value = 9.763396412722928
$ws.Cells.Item(25, 6).Value = 18.34778573295697
$ws.Cells.Item(25, 7).Value = 3.633502687565202
$ws.Cells.Item(25, 9).Value = 23.80155320932546
$ws.Cells.Item(25, 13).Value = 16.92085988271535
$ws.Cells.Item(25, 14).Value = 17.86903778332354
